# Auto-generated Excel COM-interop script
# Applies numeric updates to currentAveragePrice / Leve price & profit columns
# across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets per the Typhon_Profits dataset refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19 (hunk 1)
$ws.Range("H19").Value = 808.63635
$ws.Range("I19").Value = 688
$ws.Range("J19").Value = 909.1667
$ws.Range("K19").Value = 688
$ws.Range("L19").Value = 909.1667
$ws.Range("M19").Value = -513
$ws.Range("N19").Value = -1259.1667

# Row 33 (hunk 2)
$ws.Range("H33").Value = 235.25
$ws.Range("I33").Value = 235.25
$ws.Range("K33").Value = 235.25
$ws.Range("M33").Value = -6.25

# Row 51 (hunk 3)
$ws.Range("H51").Value = 11450
$ws.Range("I51").Value = 11450
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 11450
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -10966
$ws.Range("N51").ClearContents()

# Row 53 (hunk 4)
$ws.Range("H53").Value = 2984.7273
$ws.Range("I53").Value = 190
$ws.Range("J53").Value = 5313.6665
$ws.Range("K53").Value = 190
$ws.Range("L53").Value = 5313.6665
$ws.Range("M53").Value = 447
$ws.Range("N53").Value = -6587.6665

# Row 64 (hunk 5)
$ws.Range("H64").Value = 4928.5713
$ws.Range("I64").Value = 4500
$ws.Range("J64").Value = 4961.5386
$ws.Range("K64").Value = 4500
$ws.Range("L64").Value = 4961.5386
$ws.Range("M64").Value = -4252
$ws.Range("N64").Value = -5457.5386

# Row 67 (hunk 6)
$ws.Range("H67").Value = 4928.5713
$ws.Range("I67").Value = 4500
$ws.Range("J67").Value = 4961.5386
$ws.Range("K67").Value = 4500
$ws.Range("L67").Value = 4961.5386
$ws.Range("M67").Value = -3642
$ws.Range("N67").Value = -6677.5386

# Row 113 (hunk 7)
$ws.Range("H113").Value = 76926696
$ws.Range("I113").Value = 111113610
$ws.Range("J113").Value = 6125
$ws.Range("K113").Value = 111113610
$ws.Range("L113").Value = 6125
$ws.Range("M113").Value = -111110356
$ws.Range("N113").Value = -12633

# Row 132 (hunk 8)
$ws.Range("H132").Value = 4512.316
$ws.Range("I132").Value = 4922.2666
$ws.Range("J132").Value = 2975
$ws.Range("K132").Value = 14766.7998
$ws.Range("L132").Value = 8925
$ws.Range("M132").Value = -12236.7998
$ws.Range("N132").Value = -13985

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (hunk 9)
$ws.Range("H32").Value = 2080.4355
$ws.Range("I32").Value = 1754.0339
$ws.Range("J32").Value = 8499.666999999999
$ws.Range("K32").Value = 1754.0339
$ws.Range("L32").Value = 8499.666999999999
$ws.Range("M32").Value = -1467.0339
$ws.Range("N32").Value = -9073.666999999999

# Row 97 (hunk 10)
$ws.Range("H97").Value = 2081.25
$ws.Range("I97").Value = 1964.2858
$ws.Range("K97").Value = 1964.2858
$ws.Range("M97").Value = -1468.2858

$ws = $wb.Worksheets.Item("BSM")
# Row 80 (hunk 11)
$ws.Range("H80").Value = 638.36365
$ws.Range("I80").Value = 843.3333
$ws.Range("J80").Value = 496.46155
$ws.Range("K80").Value = 843.3333
$ws.Range("L80").Value = 496.46155
$ws.Range("M80").Value = 154.6667
$ws.Range("N80").Value = -2492.46155

# Row 83 (hunk 12)
$ws.Range("H83").Value = 638.36365
$ws.Range("I83").Value = 843.3333
$ws.Range("J83").Value = 496.46155
$ws.Range("K83").Value = 4216.6665
$ws.Range("L83").Value = 2482.30775
$ws.Range("M83").Value = 775.3334999999997
$ws.Range("N83").Value = -12466.30775

# Row 107 (hunk 13)
$ws.Range("H107").Value = 5000
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 5000
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 5000
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -8840

$ws = $wb.Worksheets.Item("CRP")
# Row 16 (hunk 14)
$ws.Range("H16").Value = 900
$ws.Range("I16").Value = 900
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 900
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -613
$ws.Range("N16").ClearContents()

# Row 31 (hunk 15)
$ws.Range("H31").Value = 2621.52
$ws.Range("I31").Value = 1389.1428
$ws.Range("J31").Value = 4190
$ws.Range("K31").Value = 1389.1428
$ws.Range("L31").Value = 4190
$ws.Range("M31").Value = -1094.1428
$ws.Range("N31").Value = -4780

# Row 34 (hunk 16)
$ws.Range("H34").Value = 2621.52
$ws.Range("I34").Value = 1389.1428
$ws.Range("J34").Value = 4190
$ws.Range("K34").Value = 1389.1428
$ws.Range("L34").Value = 4190
$ws.Range("M34").Value = -1187.1428
$ws.Range("N34").Value = -4594

# Row 62 (hunk 17)
$ws.Range("H62").Value = 3268.1
$ws.Range("I62").Value = 3216.8
$ws.Range("K62").Value = 3216.8
$ws.Range("M62").Value = -2592.8

# Row 65 (hunk 18)
$ws.Range("H65").Value = 3268.1
$ws.Range("I65").Value = 3216.8
$ws.Range("K65").Value = 16084
$ws.Range("M65").Value = -12964

# Row 88 (hunk 19)
$ws.Range("H88").Value = 35780.668
$ws.Range("J88").Value = 35780.668
$ws.Range("L88").Value = 35780.668
$ws.Range("N88").Value = -36592.668

# Row 91 (hunk 20)
$ws.Range("H91").Value = 35780.668
$ws.Range("J91").Value = 35780.668
$ws.Range("L91").Value = 35780.668
$ws.Range("N91").Value = -38588.668

# Row 107 (hunk 21)
$ws.Range("H107").Value = 1109.2
$ws.Range("I107").Value = 882.3333
$ws.Range("J107").Value = 1449.5
$ws.Range("K107").Value = 882.3333
$ws.Range("L107").Value = 1449.5
$ws.Range("M107").Value = 1037.6667
$ws.Range("N107").Value = -5289.5

# Row 113 (hunk 22)
$ws.Range("H113").Value = 900
$ws.Range("I113").Value = 900
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 900
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1270
$ws.Range("N113").ClearContents()

# Row 132 (hunk 23)
$ws.Range("H132").Value = 3324.8667
$ws.Range("I132").Value = 1287.3
$ws.Range("J132").Value = 7400
$ws.Range("K132").Value = 3861.9
$ws.Range("L132").Value = 22200
$ws.Range("M132").Value = -1331.9
$ws.Range("N132").Value = -27260

$ws = $wb.Worksheets.Item("CUL")
# Row 13 (hunk 24)
$ws.Range("H13").Value = 250
$ws.Range("I13").Value = 200
$ws.Range("J13").Value = 400
$ws.Range("K13").Value = 600
$ws.Range("L13").Value = 1200
$ws.Range("M13").Value = -432
$ws.Range("N13").Value = -1536

# Row 68 (hunk 25)
$ws.Range("H68").Value = 20400.6
$ws.Range("I68").Value = 500
$ws.Range("J68").Value = 25375.75
$ws.Range("K68").Value = 1500
$ws.Range("L68").Value = 76127.25
$ws.Range("M68").Value = -689
$ws.Range("N68").Value = -77749.25

# Row 71 (hunk 26)
$ws.Range("H71").Value = 20400.6
$ws.Range("I71").Value = 500
$ws.Range("J71").Value = 25375.75
$ws.Range("K71").Value = 4500
$ws.Range("L71").Value = 228381.75
$ws.Range("M71").Value = -444
$ws.Range("N71").Value = -236493.75

# Row 97 (hunk 27)
$ws.Range("H97").Value = 1079.0834
$ws.Range("J97").Value = 1539.8
$ws.Range("L97").Value = 4619.4
$ws.Range("N97").Value = -5611.4

# Row 112 (hunk 28)
$ws.Range("H112").Value = 1625
$ws.Range("I112").Value = 1200
$ws.Range("J112").Value = 2900
$ws.Range("K112").Value = 3600
$ws.Range("L112").Value = 8700
$ws.Range("M112").Value = -2492
$ws.Range("N112").Value = -10916

# Row 121 (hunk 29)
$ws.Range("H121").Value = 5231.091
$ws.Range("I121").Value = 592.5
$ws.Range("J121").Value = 5694.95
$ws.Range("K121").Value = 1777.5
$ws.Range("L121").Value = 17084.85
$ws.Range("M121").Value = -467.5
$ws.Range("N121").Value = -19704.85

# Row 131 (hunk 30)
$ws.Range("H131").Value = 785.1900000000001
$ws.Range("J131").Value = 786.0505000000001
$ws.Range("L131").Value = 2358.1515
$ws.Range("N131").Value = -12438.1515

$ws = $wb.Worksheets.Item("GSM")
# Row 102 (hunk 31)
$ws.Range("H102").Value = 1748.6207
$ws.Range("I102").Value = 1738.8462
$ws.Range("K102").Value = 1738.8462
$ws.Range("M102").Value = -116.8462

$ws = $wb.Worksheets.Item("LTW")
# Row 22 (hunk 32)
$ws.Range("H22").Value = 2423.4666
$ws.Range("I22").Value = 1905.1
$ws.Range("J22").Value = 3460.2
$ws.Range("K22").Value = 1905.1
$ws.Range("L22").Value = 3460.2
$ws.Range("M22").Value = -1610.1
$ws.Range("N22").Value = -4050.2

# Row 27 (hunk 33)
$ws.Range("H27").Value = 2423.4666
$ws.Range("I27").Value = 1905.1
$ws.Range("J27").Value = 3460.2
$ws.Range("K27").Value = 1905.1
$ws.Range("L27").Value = 3460.2
$ws.Range("M27").Value = -1798.1
$ws.Range("N27").Value = -3674.2

# Row 40 (hunk 34)
$ws.Range("H40").Value = 3609.5625
$ws.Range("I40").Value = 2791.5
$ws.Range("J40").Value = 4100.4
$ws.Range("K40").Value = 2791.5
$ws.Range("L40").Value = 4100.4
$ws.Range("M40").Value = -2655.5
$ws.Range("N40").Value = -4372.4

# Row 61 (hunk 35)
$ws.Range("H61").Value = 5657.857
$ws.Range("I61").Value = 3150
$ws.Range("J61").Value = 9001.666999999999
$ws.Range("K61").Value = 3150
$ws.Range("L61").Value = 9001.666999999999
$ws.Range("M61").Value = -2948
$ws.Range("N61").Value = -9405.666999999999

# Row 111 (hunk 36)
$ws.Range("H111").Value = 32387
$ws.Range("J111").Value = 32387
$ws.Range("L111").Value = 32387
$ws.Range("N111").Value = -40567

# Row 113 (hunk 37)
$ws.Range("H113").Value = 5657.857
$ws.Range("I113").Value = 3150
$ws.Range("J113").Value = 9001.666999999999
$ws.Range("K113").Value = 3150
$ws.Range("L113").Value = 9001.666999999999
$ws.Range("M113").Value = -980
$ws.Range("N113").Value = -13341.667

$ws = $wb.Worksheets.Item("WVR")
# Row 132 (hunk 38)
$ws.Range("H132").Value = 1371.8235
$ws.Range("I132").Value = 801.75
$ws.Range("J132").Value = 1878.5555
$ws.Range("K132").Value = 2405.25
$ws.Range("L132").Value = 5635.666499999999
$ws.Range("M132").Value = 124.75
$ws.Range("N132").Value = -10695.6665

Write-Output "Applied Typhon_Profits numeric updates across 38 rows."